$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("F2").Value = "-"

$ws.Range("E4").Value = "['MCT-3A-Processos de Usinagem 1', -, -, -]"

$ws.Range("E6").Value = "['MCT-3A-Processos de Usinagem 1', -, -, -]"

$ws.Range("E7").Value = "['MCT-3A-Processos de Usinagem 1', -, -, -]"

$ws.Range("B8").Value = "-"
$ws.Range("E8").Value = "['MCT-3A-Processos de Usinagem 1', -, -, -]"
$ws.Range("F8").Value = "-"
